{"js": "// The document's single table holds a 20x5 grid of arithmetic prompts\n// (\"48+25=\", \"60-12=\", ...). The commit replaces every prompt's text with\n// a new expression, in row-major (left-to-right, top-to-bottom) cell order.\n// NEW_VALUES holds those 100 replacement strings in that same order.\nconst NEW_VALUES = [\n  \"19+2=\", \"88-39=\", \"37+8=\", \"7+87=\", \"17+44=\", \"76+16=\", \"90-81=\", \"92-86=\", \"54+29=\", \"71-8=\",\n  \"84-68=\", \"46+47=\", \"68+24=\", \"91-18=\", \"27+15=\", \"49+19=\", \"87+4=\", \"63-14=\", \"55-28=\", \"41-38=\",\n  \"87+6=\", \"22-14=\", \"73-67=\", \"45-19=\", \"90-22=\", \"9+74=\", \"75-49=\", \"75-7=\", \"24+59=\", \"92-8=\",\n  \"7+56=\", \"36-19=\", \"65-57=\", \"19+3=\", \"38-9=\", \"8+68=\", \"41-18=\", \"5+37=\", \"46+38=\", \"81-57=\",\n  \"39+55=\", \"71-49=\", \"49+32=\", \"28+55=\", \"84+8=\", \"15+17=\", \"61-48=\", \"29+9=\", \"65-48=\", \"8+57=\",\n  \"50-12=\", \"84-37=\", \"61-13=\", \"71-54=\", \"37+19=\", \"42-16=\", \"9+15=\", \"93-57=\", \"67+6=\", \"6+8=\",\n  \"5+9=\", \"39+6=\", \"74-26=\", \"90-18=\", \"39+39=\", \"71-6=\", \"4+58=\", \"39+38=\", \"15+37=\", \"13+48=\",\n  \"57-9=\", \"77+14=\", \"19+29=\", \"4+38=\", \"15+29=\", \"46+47=\", \"16+16=\", \"80-26=\", \"55-46=\", \"39+35=\",\n  \"34-27=\", \"56+26=\", \"14+37=\", \"40-25=\", \"90-86=\", \"92-43=\", \"17+24=\", \"87+8=\", \"36+8=\", \"63-45=\",\n  \"68+16=\", \"46+8=\", \"6+69=\", \"81-47=\", \"46-29=\", \"60-42=\", \"84-15=\", \"24-17=\", \"31-14=\", \"28+68=\"\n];\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nlet idx = 0;\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    if (idx < NEW_VALUES.length) {\n      // TableCell.value overwrites the cell's text while keeping the\n      // existing run/paragraph formatting (font, size, alignment, etc).\n      cell.value = NEW_VALUES[idx];\n    }\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# The document's single table holds a 20x5 grid of arithmetic prompts\n# (\"48+25=\", \"60-12=\", ...). The commit replaces every prompt's text with\n# a new expression, in row-major (left-to-right, top-to-bottom) cell order.\n# $NewValues holds those 100 replacement strings in that same order.\n$NewValues = @(\n    \"19+2=\",\"88-39=\",\"37+8=\",\"7+87=\",\"17+44=\",\"76+16=\",\"90-81=\",\"92-86=\",\n    \"54+29=\",\"71-8=\",\"84-68=\",\"46+47=\",\"68+24=\",\"91-18=\",\"27+15=\",\"49+19=\",\n    \"87+4=\",\"63-14=\",\"55-28=\",\"41-38=\",\"87+6=\",\"22-14=\",\"73-67=\",\"45-19=\",\n    \"90-22=\",\"9+74=\",\"75-49=\",\"75-7=\",\"24+59=\",\"92-8=\",\"7+56=\",\"36-19=\",\n    \"65-57=\",\"19+3=\",\"38-9=\",\"8+68=\",\"41-18=\",\"5+37=\",\"46+38=\",\"81-57=\",\n    \"39+55=\",\"71-49=\",\"49+32=\",\"28+55=\",\"84+8=\",\"15+17=\",\"61-48=\",\"29+9=\",\n    \"65-48=\",\"8+57=\",\"50-12=\",\"84-37=\",\"61-13=\",\"71-54=\",\"37+19=\",\"42-16=\",\n    \"9+15=\",\"93-57=\",\"67+6=\",\"6+8=\",\"5+9=\",\"39+6=\",\"74-26=\",\"90-18=\",\n    \"39+39=\",\"71-6=\",\"4+58=\",\"39+38=\",\"15+37=\",\"13+48=\",\"57-9=\",\"77+14=\",\n    \"19+29=\",\"4+38=\",\"15+29=\",\"46+47=\",\"16+16=\",\"80-26=\",\"55-46=\",\"39+35=\",\n    \"34-27=\",\"56+26=\",\"14+37=\",\"40-25=\",\"90-86=\",\"92-43=\",\"17+24=\",\"87+8=\",\n    \"36+8=\",\"63-45=\",\"68+16=\",\"46+8=\",\"6+69=\",\"81-47=\",\"46-29=\",\"60-42=\",\n    \"84-15=\",\"24-17=\",\"31-14=\",\"28+68=\"\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    if ($idx -lt $NewValues.Length) {\n      # Cell.Range.Text replaces only the cell's content (preserving the\n      # trailing cell-mark) and keeps the existing run formatting.\n      $table.Cell($r, $c).Range.Text = $NewValues[$idx]\n    }\n    $idx++\n  }\n}\n"}
